$d = $word.ActiveDocument

# The bibliography paragraph ("Apostila ou texto fornecido pelo(s)
# docente(s) ...") is kept. Everything that used to follow it on the page
# -- a blank separator paragraph, the "Ver no Jupiter Salvar em pdf Salvar
# em docx" line, and the "(c) 2020 ... Creative Commons Attribution" site
# footer line -- must be deleted, leaving the blank paragraph that precedes
# the trailing page-break paragraph untouched.

$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Apostila ou texto fornecido*") {
        $startPara = $i
        break
    }
}

$endPara = $null
for ($i = $startPara + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Creative Commons Attribution*") {
        $endPara = $i
        break
    }
}

$rangeStart = $d.Paragraphs.Item($startPara + 1).Range.Start
$rangeEnd = $d.Paragraphs.Item($endPara).Range.End
$d.Range($rangeStart, $rangeEnd).Delete()
